$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("teams")

# Headers (row 1)
$ws.Range("A1").Value = "Week"
$ws.Range("B1").Value = "Date"
$ws.Range("C1").Value = "EarlyHomeTeam"
$ws.Range("D1").Value = "EarlyHomeScore"
$ws.Range("E1").Value = "EarlyAwayTeam"
$ws.Range("F1").Value = "EarlyAwayScore"
$ws.Range("G1").Value = "MidHomeTeam"
$ws.Range("H1").Value = "MidHomeScore"
$ws.Range("I1").Value = "MidAwayTeam"
$ws.Range("J1").Value = "MidAwayScore"
$ws.Range("K1").Value = "LateHomeTeam"
$ws.Range("L1").Value = "LateHomeScore"
$ws.Range("M1").Value = "LateAwayTeam"
$ws.Range("N1").Value = "LateAwayScore"

# Row 2 data
$ws.Range("A2").Value = 1
$ws.Range("C2").Value = "Team 1"
$ws.Range("D2").Value = 76
$ws.Range("E2").Value = "Team 2"
$ws.Range("F2").Value = 70
$ws.Range("G2").Value = "Team 3"
$ws.Range("H2").Value = 90
$ws.Range("I2").Value = "Team 4"
$ws.Range("J2").Value = 94
$ws.Range("K2").Value = "Team 5"
$ws.Range("L2").Value = 55
$ws.Range("M2").Value = "Team 6"
$ws.Range("N2").Value = 50

# B2 (Date) must stay literal text "06-04-2025", not auto-converted to a date serial.
# Enter it as a formula producing the literal string, then paste-special as values
# so it collapses to a plain shared-string cell with no number format / style applied.
$ws.Range("B2").Formula = "=""06-04-2025"""
$ws.Range("B2").Copy()
$ws.Range("B2").PasteSpecial(-4163)
$excel.CutCopyMode = 0

# Column widths to autofit the new table
$ws.Range("A1:N2").Columns.AutoFit()

# Select whole sheet then fix the selection to A1 like a fresh CSV paste, and make "teams" the active/visible tab
$ws.Range("A1:XFD1048576").Select()
$ws.Activate()
